$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new literal text value, taken from the updated
# cryptos price/volume snapshot. All of these cells hold plain text (prices
# like "42.149.08" and padded percentages like "  +1.83%  " are not valid
# numeric literals in this sheet), so each is written with a leading
# apostrophe to force text entry, then the style is reset to "Normal" so no
# stray NumberFormat/quote-prefix styling is left behind on the cell.
$updates = [ordered]@{
    "D2" = '42.149.08';
    "E2" = '  +1.83%  ';
    "D3" = '2.220.57';
    "E3" = '  +1.27%  ';
    "E4" = '  -0.03%  ';
    "D5" = '250.86';
    "E5" = '  -1.40%  ';
    "D7" = '68.34';
    "E7" = '  -0.55%  ';
    "E9" = '  +6.03%  ';
    "D10" = '39.45';
    "E10" = '  +3.48%  ';
    "E11" = '  +0.96%  ';
    "D12" = '0.0940';
    "E12" = '  -0.07%  ';
    "D13" = '7.08';
    "E13" = '  -1.42%  ';
    "D14" = '0.104';
    "E14" = '  -0.95%  ';
    "D15" = '2.555.19';
    "E15" = '  +1.53%  ';
    "D16" = '0.872';
    "E16" = '  -0.43%  ';
    "D17" = '14.52';
    "E17" = '  -0.62%  ';
    "D18" = '2.215.03';
    "E18" = '  +1.48%  ';
    "D19" = '42.033.77';
    "E19" = '  +1.76%  ';
    "E20" = '  +0.46%  ';
    "D21" = '72.39';
    "E21" = '  +0.14%  ';
    "E22" = '  -1.59%  ';
    "D23" = '232.07';
    "E23" = '  -0.54%  ';
    "D24" = '2.04';
    "E24" = '  -0.45%  ';
    "E25" = '  -0.48%  ';
    "E26" = '  +0.02%  ';
    "E27" = '  -5.71%  ';
    "E28" = '  -4.72%  ';
    "E29" = '  -1.88%  ';
    "E30" = '  -2.05%  ';
    "D31" = '166.95';
    "E31" = '  -1.81%  ';
    "D32" = '20.47';
    "E32" = '  -1.15%  ';
    "D33" = '6.16';
    "E33" = '  +10.74%  ';
    "D34" = '0.121';
    "E34" = '  +0.93%  ';
    "D35" = '0.0782';
    "E35" = '  +5.98%  ';
    "E36" = '  -1.04%  ';
    "D37" = '26.82';
    "E37" = '  +0.18%  ';
    "D38" = '4.62';
    "E38" = '  -0.37%  ';
    "D39" = '4.12';
    "E39" = '  -0.02%  ';
    "E40" = '  +4.07%  ';
    "D41" = '2.25';
    "E41" = '  +1.48%  ';
    "D42" = '12.29';
    "E42" = '  -4.27%  ';
    "D43" = '5.70';
    "E43" = '  -0.52%  ';
    "E44" = '  +1.15%  ';
    "D45" = '62.01';
    "E45" = '  -3.98%  ';
    "E46" = '  -3.66%  ';
    "E47" = '  -0.73%  ';
    "E48" = '  -0.46%  ';
    "E49" = '  -0.21%  ';
    "E50" = '  +1.50%  ';
    "B51" = 'TrustWalletToken';
    "C51" = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt';
    "D51" = '1.17';
    "E51" = '  -1.00%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.Value = "'" + $updates[$ref]
    $cell.Style = "Normal"
}
